$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.937888198757764
$ws.Range("C2").Value = 0.934010152284264
$ws.Range("D2").Value = 0.935064935064935
$ws.Range("E2").Value = 0.882352941176471
